$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record needs to be inserted before the existing row 124.
# Insert a blank row at row 124, shifting rows 124:143 down to 125:144.
$ws.Rows(124).Insert()

# Populate the newly inserted row 124 with the new "Puerro" observation.
$ws.Range("A124").Value = 10
$ws.Range("B124").Value = "Vega Modelo de Temuco"
$ws.Range("C124").Value = "La Araucanía"
$ws.Range("D124").Value = 44491
$ws.Range("E124").Value = 9
$ws.Range("F124").Value = 100112005
$ws.Range("G124").Value = "Puerro"
$ws.Range("H124").Value = "Azul de Maquehue"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 65
$ws.Range("K124").Value = 7000
$ws.Range("L124").Value = 7000
$ws.Range("M124").Value = 7000
$ws.Range("N124").Value = "$/docena de paquetes"
$ws.Range("O124").Value = "Provincia de Cautín"
$ws.Range("P124").Value = 583
$ws.Range("Q124").Value = 12
$ws.Range("R124").Value = "Hortaliza"
